$d = $word.ActiveDocument

$old = "Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός του Περσέα 2022: 16-25 Ιανουαρίου, 7-16 Νοεμβρίου, 6-15 Δεκεμβρίου"
$new = "2022 Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός του Περσέα: 16-25 Ιανουαρίου, 7-16 Νοεμβρίου, 6-15 Δεκεμβρίου"

$d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
